function Set-TextCell {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws "D2" "57.442.22"
Set-TextCell $ws "E2" "  -0.73%  "
Set-TextCell $ws "D3" "3.090.79"
Set-TextCell $ws "E3" "  +0.94%  "
Set-TextCell $ws "E4" "  -0.01%  "
Set-TextCell $ws "D5" "521.36"
Set-TextCell $ws "E5" "  +0.81%  "
Set-TextCell $ws "D6" "140.75"
Set-TextCell $ws "E6" "  -1.29%  "
Set-TextCell $ws "E7" "  -0.01%  "
Set-TextCell $ws "D8" "3.089.35"
Set-TextCell $ws "E8" "  +0.96%  "
Set-TextCell $ws "E9" "  +0.20%  "
Set-TextCell $ws "E10" "  -1.56%  "
Set-TextCell $ws "E11" "  -0.70%  "
Set-TextCell $ws "E12" "  +1.72%  "
Set-TextCell $ws "D13" "3.620.22"
Set-TextCell $ws "E13" "  +0.79%  "
Set-TextCell $ws "E14" "  +0.98%  "
Set-TextCell $ws "E15" "  -2.07%  "
Set-TextCell $ws "D16" "0.0000163"
Set-TextCell $ws "E16" "  -0.63%  "
Set-TextCell $ws "D17" "57.519.20"
Set-TextCell $ws "E17" "  -0.64%  "
Set-TextCell $ws "D18" "3.089.27"
Set-TextCell $ws "E18" "  +0.33%  "
Set-TextCell $ws "E19" "  -0.63%  "
Set-TextCell $ws "D20" "12.72"
Set-TextCell $ws "E20" "  -1.00%  "
Set-TextCell $ws "D21" "8.03"
Set-TextCell $ws "E21" "  -0.75%  "
Set-TextCell $ws "D22" "338.29"
Set-TextCell $ws "E22" "  +1.81%  "
Set-TextCell $ws "E23" "  -0.06%  "
Set-TextCell $ws "E24" "  +1.80%  "
Set-TextCell $ws "E25" "  +1.71%  "
Set-TextCell $ws "E26" "  -1.70%  "
Set-TextCell $ws "E27" "  +0.31%  "
Set-TextCell $ws "D28" "0.0₃0907"
Set-TextCell $ws "E28" "  +0.45%  "
Set-TextCell $ws "E29" "  -0.02%  "
Set-TextCell $ws "D30" "6.45"
Set-TextCell $ws "E30" "  -0.32%  "
Set-TextCell $ws "E31" "  -1.72%  "
Set-TextCell $ws "E32" "  +1.69%  "
Set-TextCell $ws "D33" "20.82"
Set-TextCell $ws "E33" "  +0.56%  "
Set-TextCell $ws "E34" "  -1.66%  "
Set-TextCell $ws "D35" "155.99"
Set-TextCell $ws "E35" "  +0.85%  "
Set-TextCell $ws "D36" "4.60"
Set-TextCell $ws "E36" "  +1.22%  "
Set-TextCell $ws "E37" "  +1.35%  "
Set-TextCell $ws "D38" "27.04"
Set-TextCell $ws "E38" "  -0.08%  "
Set-TextCell $ws "D39" "1.24"
Set-TextCell $ws "E39" "  -2.32%  "
Set-TextCell $ws "D40" "0.0659"
Set-TextCell $ws "E40" "  -3.06%  "
Set-TextCell $ws "E41" "  +10.49%  "
Set-TextCell $ws "B42" "Filecoin"
Set-TextCell $ws "C42" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D42" "3.93"
Set-TextCell $ws "E42" "  -0.07%  "
Set-TextCell $ws "B43" "RenzoRestakedETH"
Set-TextCell $ws "C43" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextCell $ws "D43" "3.131.15"
Set-TextCell $ws "E43" "  +0.83%  "
Set-TextCell $ws "E44" "  +3.81%  "
Set-TextCell $ws "E45" "  +0.59%  "
Set-TextCell $ws "D46" "0.999"
Set-TextCell $ws "E46" "  -0.06%  "
Set-TextCell $ws "D47" "2.297.28"
Set-TextCell $ws "E47" "  +1.01%  "
Set-TextCell $ws "D48" "0.0257"
Set-TextCell $ws "E48" "  -0.75%  "
Set-TextCell $ws "E49" "  +3.45%  "
Set-TextCell $ws "D50" "20.46"
Set-TextCell $ws "E50" "  -1.68%  "
Set-TextCell $ws "E51" "  +1.04%  "
